$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to remain plain text (matches original inlineStr cells)
# by temporarily marking the cell Text-formatted, assigning the literal value,
# then resetting the style back to Normal so no stray number-format/style survives.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '56.965.19'
Set-TextValue $ws.Range("E2") '  +0.42%  '

Set-TextValue $ws.Range("D3") '2.983.46'
Set-TextValue $ws.Range("E3") '  -0.55%  '

Set-TextValue $ws.Range("E4") '  +0.25%  '

Set-TextValue $ws.Range("D5") '500.79'
Set-TextValue $ws.Range("E5") '  -1.77%  '

Set-TextValue $ws.Range("D6") '137.95'
Set-TextValue $ws.Range("E6") '  -0.45%  '

Set-TextValue $ws.Range("E7") '  +0.13%  '

Set-TextValue $ws.Range("D8") '0.431'
Set-TextValue $ws.Range("E8") '  -0.84%  '

Set-TextValue $ws.Range("D9") '7.49'
Set-TextValue $ws.Range("E9") '  +0.38%  '

Set-TextValue $ws.Range("D10") '0.108'
Set-TextValue $ws.Range("E10") '  +0.38%  '

Set-TextValue $ws.Range("D11") '0.360'
Set-TextValue $ws.Range("E11") '  +1.22%  '

Set-TextValue $ws.Range("D12") '3.510.80'
Set-TextValue $ws.Range("E12") '  -0.06%  '

Set-TextValue $ws.Range("E13") '  -0.84%  '

Set-TextValue $ws.Range("D14") '25.94'
Set-TextValue $ws.Range("E14") '  +0.63%  '

Set-TextValue $ws.Range("D15") '0.0000159'
Set-TextValue $ws.Range("E15") '  +1.83%  '

Set-TextValue $ws.Range("D16") '57.107.23'
Set-TextValue $ws.Range("E16") '  +0.59%  '

Set-TextValue $ws.Range("D17") '6.10'
Set-TextValue $ws.Range("E17") '  +2.67%  '

Set-TextValue $ws.Range("D18") '2.999.99'
Set-TextValue $ws.Range("E18") '  +0.05%  '

Set-TextValue $ws.Range("D19") '12.65'
Set-TextValue $ws.Range("E19") '  +1.36%  '

Set-TextValue $ws.Range("D20") '7.81'
Set-TextValue $ws.Range("E20") '  -0.14%  '

Set-TextValue $ws.Range("D21") '322.80'
Set-TextValue $ws.Range("E21") '  -0.94%  '

Set-TextValue $ws.Range("D22") '0.997'
Set-TextValue $ws.Range("E22") '  -0.25%  '

Set-TextValue $ws.Range("D23") '5.65'
Set-TextValue $ws.Range("E23") '  -0.53%  '

Set-TextValue $ws.Range("D24") '0.490'
Set-TextValue $ws.Range("E24") '  +1.05%  '

Set-TextValue $ws.Range("D25") '64.02'
Set-TextValue $ws.Range("E25") '  +1.25%  '

Set-TextValue $ws.Range("B26") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D26") '1.01'
Set-TextValue $ws.Range("E26") '  +1.30%  '

Set-TextValue $ws.Range("B27") 'Kaspa'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D27") '0.163'
Set-TextValue $ws.Range("E27") '  -4.07%  '

Set-TextValue $ws.Range("D28") '0.0₃0900'
Set-TextValue $ws.Range("E28") '  -0.38%  '

Set-TextValue $ws.Range("D29") '6.62'
Set-TextValue $ws.Range("E29") '  -0.29%  '

Set-TextValue $ws.Range("D30") '7.15'
Set-TextValue $ws.Range("E30") '  +2.32%  '

Set-TextValue $ws.Range("D31") '1.78'
Set-TextValue $ws.Range("E31") '  -1.17%  '

Set-TextValue $ws.Range("D32") '1.17'
Set-TextValue $ws.Range("E32") '  -4.01%  '

Set-TextValue $ws.Range("D33") '20.23'
Set-TextValue $ws.Range("E33") '  -1.98%  '

Set-TextValue $ws.Range("D34") '4.66'
Set-TextValue $ws.Range("E34") '  +2.60%  '

Set-TextValue $ws.Range("D35") '153.64'
Set-TextValue $ws.Range("E35") '  -0.42%  '

Set-TextValue $ws.Range("D36") '5.79'
Set-TextValue $ws.Range("E36") '  +2.08%  '

Set-TextValue $ws.Range("D37") '1.25'
Set-TextValue $ws.Range("E37") '  -1.12%  '

Set-TextValue $ws.Range("D38") '24.20'
Set-TextValue $ws.Range("E38") '  +1.90%  '

Set-TextValue $ws.Range("D39") '0.0668'
Set-TextValue $ws.Range("E39") '  -1.60%  '

Set-TextValue $ws.Range("D40") '3.017.76'
Set-TextValue $ws.Range("E40") '  -0.54%  '

Set-TextValue $ws.Range("D41") '37.56'
Set-TextValue $ws.Range("E41") '  +1.38%  '

Set-TextValue $ws.Range("E42") '  +0.23%  '

Set-TextValue $ws.Range("D43") '3.78'
Set-TextValue $ws.Range("E43") '  +2.66%  '

Set-TextValue $ws.Range("D44") '0.643'
Set-TextValue $ws.Range("E44") '  -0.72%  '

Set-TextValue $ws.Range("D45") '2.212.88'
Set-TextValue $ws.Range("E45") '  -3.02%  '

Set-TextValue $ws.Range("D46") '1.39'
Set-TextValue $ws.Range("E46") '  -1.71%  '

Set-TextValue $ws.Range("D47") '0.953'
Set-TextValue $ws.Range("E47") '  -4.55%  '

Set-TextValue $ws.Range("D48") '5.96'
Set-TextValue $ws.Range("E48") '  +1.66%  '

Set-TextValue $ws.Range("D49") '0.0235'
Set-TextValue $ws.Range("E49") '  -1.14%  '

Set-TextValue $ws.Range("D50") '19.21'
Set-TextValue $ws.Range("E50") '  -0.06%  '

Set-TextValue $ws.Range("D51") '1.82'
Set-TextValue $ws.Range("E51") '  -6.67%  '
